# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
# The "Periodo Mora" column (E16:E19) previously listed periods in
# descending order (2506, 2505, 2504, 2503). Update it so the periods
# read in ascending order (2503, 2504, 2505, 2506).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "2503"
$ws.Range("E17").Value = "2504"
$ws.Range("E18").Value = "2505"
$ws.Range("E19").Value = "2506"
